$wb = $excel.ActiveWorkbook

# Antibody sheet: insert new column "antibodySpecificity" before the
# existing "targetedAntigen" column (D), shifting D..AZ to E..BA.
$wsAntibody = $wb.Worksheets.Item("Antibody")
$wsAntibody.Range("D1").EntireColumn.Insert()
$wsAntibody.Range("D1").Value = "antibodySpecificity"

# Hybridoma sheet: insert new column "antibodySpecificity" before the
# existing "targetedAntigen" column (E), shifting E..BA to F..BB.
$wsHybridoma = $wb.Worksheets.Item("Hybridoma")
$wsHybridoma.Range("E1").EntireColumn.Insert()
$wsHybridoma.Range("E1").Value = "antibodySpecificity"
